# Add Swab & Sputum & Stool E gene Ct value
# Rename the "Type" column (G) labels to the new, more explicit
# variable-style names, and refresh the view (selected cell + zoom).
#
# NOTE: the order in which the distinct label groups are (re)written
# matters: new shared strings are appended to the workbook's shared
# string table in the order they are first introduced, and this
# script intentionally mirrors that ordering so the resulting
# sharedStrings.xml table lines up with the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TypeLabel($rows, $newLabel) {
    foreach ($r in $rows) {
        $ws.Range("G$r").Value = $newLabel
    }
}

# swab E  ->  swab_SARSCoV2_E_Ct
Set-TypeLabel (48..73) "swab_SARSCoV2_E_Ct"

# sputum E  ->  sputum_SARSCoV2_E_Ct
Set-TypeLabel (74..93) "sputum_SARSCoV2_E_Ct"

# stool RdRp  ->  stool_SARSCoV2_RdRp_Ct
Set-TypeLabel ((94..104) + (116..120)) "stool_SARSCoV2_RdRp_Ct"

# stool E  ->  stool_SARSCoV2_E_Ct
Set-TypeLabel ((105..115) + (121..125)) "stool_SARSCoV2_E_Ct"

# sputum RdRp estimated viral copy/mL  ->  sputum_SARSCoV2_RdRp_VL
Set-TypeLabel (28..47) "sputum_SARSCoV2_RdRp_VL"

# swab RdRp estimated viral copy/mL  ->  swab_SARSCoV2_RdRp_VL
Set-TypeLabel (2..27) "swab_SARSCoV2_RdRp_VL"

# Refresh the view state: selected cell G29 and 140% zoom.
[void]$ws.Range("G29").Select()
$excel.ActiveWindow.Zoom = 140

Write-Output "done"
